# 每日学习.xlsx - "Add files via upload" commit
# Adds three new daily-log rows (34-36) to Sheet1, tweaks the text of an
# existing log entry (C32), and moves the active selection to the new
# last-used cell (E36), mirroring how Excel leaves the view after the
# author typed the new rows in and saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Small wording tweak to an already-existing cell -----------------
# "到了微分的定义" (reached the definition of differentiation)
# becomes "到了求导法则" (reached the rules of differentiation)
$ws.Range("C32").Value = "到了求导法则"

# --- New row 34: first entry of the new day ---------------------------
$ws.Range("A34").Value = 1
$ws.Range("B34").Value = "有点摸鱼，看了两篇英语短文"
$ws.Range("C34").Value = "背了小半小时的单词"

# --- New row 35 ---------------------------------------------------------
$ws.Range("A35").Value = 2
$ws.Range("B35").Value = "3：07 - 6：09"
$ws.Range("C35").Value = "数学：求导那边 "
$ws.Range("E35").Value = "书： 看到了1.1.7(不知道为啥，感觉比C语言那本书难看很多，看不大明白)"

# --- New row 36 ---------------------------------------------------------
$ws.Range("A36").Value = 3
$ws.Range("B36").Value = "2：59 —5：37"
$ws.Range("C36").Value = "微分刚开始"
$ws.Range("E36").Value = "书的话，不出意外每天两页，我就不说了"

# --- Leave the selection where the author's cursor ended up -----------
[void]$ws.Range("E36").Select()
